# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Row 3 corresponds to the f1cad7fc-... file in all three sheets.
# Status switches from "Ready for handoff" to "Handback transform failed"
# because the handback transform for that file failed.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Error Detail column (P) gets a new diagnostic message for the failed handback.
$zhcn.Range("P3").Value = "Handback file name: k1sdnnq0.ff0 is different with handoff file name: f1cad7fc-42a4-4e7a-ae58-241869796338.524c9a4474c973c5a2bf69d41bc1dac29e5ab262.zh-cn."
$dede.Range("P3").Value = "Handback file name: k1sdnnq0.ff0 is different with handoff file name: f1cad7fc-42a4-4e7a-ae58-241869796338.524c9a4474c973c5a2bf69d41bc1dac29e5ab262.de-de."

# Widen the Error Detail column (P) on both locale sheets to fit the new text.
# (39.15 is the COM ColumnWidth input that round-trips to a stored width of
# exactly 40 in the saved OOXML, matching Excel's own column-width autofit.)
$zhcn.Range("P1:P3").ColumnWidth = 39.15
$dede.Range("P1:P3").ColumnWidth = 39.15
